$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.56
$ws.Range("F2").Value = 0.5596767676767677
$ws.Range("G2").Value = 0.56
$ws.Range("H2").Value = 0.5589401846647932

$ws.Range("E3").Value = 0.52
$ws.Range("F3").Value = 0.5208992372541148
$ws.Range("G3").Value = 0.52
$ws.Range("H3").Value = 0.5196158463385354

$ws.Range("E4").Value = 0.59
$ws.Range("F4").Value = 0.590632183908046
$ws.Range("G4").Value = 0.59
$ws.Range("H4").Value = 0.5873959068454482

$ws.Range("E5").Value = 0.51
$ws.Range("F5").Value = 0.5121428571428571
$ws.Range("G5").Value = 0.51
$ws.Range("H5").Value = 0.5068977992161592

$ws.Range("E7").Value = 0.6
$ws.Range("F7").Value = 0.6004801920768308
$ws.Range("G7").Value = 0.6
$ws.Range("H7").Value = 0.6

$ws.Range("E8").Value = 0.66
$ws.Range("F8").Value = 0.6476984126984127
$ws.Range("G8").Value = 0.66
$ws.Range("H8").Value = 0.6469289914066034

$ws.Range("E9").Value = 0.64
$ws.Range("F9").Value = 0.6255158730158731
$ws.Range("G9").Value = 0.64
$ws.Range("H9").Value = 0.6261601085481683

$ws.Range("E13").Value = 0.66
$ws.Range("F13").Value = 0.6458275058275059
$ws.Range("G13").Value = 0.66
$ws.Range("H13").Value = 0.6340952380952382

$ws.Range("E14").Value = 0.56
$ws.Range("F14").Value = 0.5640465676435167
$ws.Range("G14").Value = 0.56
$ws.Range("H14").Value = 0.56

$ws.Range("E15").Value = 0.55
$ws.Range("F15").Value = 0.5518
$ws.Range("G15").Value = 0.55
$ws.Range("H15").Value = 0.5504053648283456

$ws.Range("E18").Value = 0.57
$ws.Range("F18").Value = 0.5819166666666666
$ws.Range("G18").Value = 0.57
$ws.Range("H18").Value = 0.5660677318862426

$ws.Range("E19").Value = 0.63
$ws.Range("F19").Value = 0.6356924315619967
$ws.Range("G19").Value = 0.63
$ws.Range("H19").Value = 0.6297409740974098

$ws.Range("F20").Value = 0.7228368794326241
$ws.Range("H20").Value = 0.6731112686038272

$ws.Range("E21").Value = 0.71
$ws.Range("F21").Value = 0.6727843137254902
$ws.Range("G21").Value = 0.71
$ws.Range("H21").Value = 0.6781691601244261

$ws.Range("E22").Value = 0.7
$ws.Range("F22").Value = 0.655813953488372
$ws.Range("G22").Value = 0.7
$ws.Range("H22").Value = 0.6632911392405063

$ws.Range("E25").Value = 0.73
$ws.Range("F25").Value = 0.7136769759450172
$ws.Range("G25").Value = 0.73
$ws.Range("H25").Value = 0.641099446459248

$ws.Range("E26").Value = 0.64
$ws.Range("F26").Value = 0.6511764705882354
$ws.Range("G26").Value = 0.64
$ws.Range("H26").Value = 0.628

$ws.Range("E27").Value = 0.67
$ws.Range("F27").Value = 0.6703875968992247
$ws.Range("G27").Value = 0.67
$ws.Range("H27").Value = 0.6685028732735154

$ws.Range("E28").Value = 0.61
$ws.Range("F28").Value = 0.6137142857142857
$ws.Range("G28").Value = 0.61
$ws.Range("H28").Value = 0.6011244979919679

$ws.Range("E31").Value = 0.58
$ws.Range("F31").Value = 0.5905197505197505
$ws.Range("G31").Value = 0.58
$ws.Range("H31").Value = 0.5542342342342342
